# Fix grammar error in PP: remove the redundant "for those values " phrase
# from the bullet on Slide 3 ("Only data that required the comparison of ABV
# to IBU had observations with the value of <NA> for those values were
# dropped ") and split the remaining sentence into three runs, matching the
# author's edit:
#   1) "Only data that required the comparison of ABV to IBU had
#       observations with the value of <"
#   2) "NA> "
#   3) "were dropped "

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(5)            # "Content Placeholder 2"
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(2)             # the "Only data that required..." bullet

# Sanity check / locate the text to remove regardless of exact offsets.
$fullText = $para.Text
$needle = "for those values "
$startPos = $fullText.IndexOf($needle)

if ($startPos -ge 0) {
    # TextRange.Characters uses 1-based character positions.
    $toDelete = $para.Characters($startPos + 1, $needle.Length)
    $toDelete.Delete()
}

# Recompute lengths for the three target runs after the deletion.
$part1 = "Only data that required the comparison of ABV to IBU had observations with the value of <"
$part2 = "NA> "
$part3 = "were dropped "

$run1 = $para.Characters(1, $part1.Length)
$run2 = $para.Characters($part1.Length + 1, $part2.Length)
$run3 = $para.Characters($part1.Length + $part2.Length + 1, $part3.Length)

$run1.Text = $part1
$run2.Text = $part2
$run3.Text = $part3
